$wb = $excel.ActiveWorkbook

# --- Section_A updates ---
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("D2").Value = "CS461"
$wsA.Range("F2").Value = "Free"

$wsA.Range("B3").Value = "CS461"
$wsA.Range("C3").Value = "CS304"
$wsA.Range("D3").Value = "Free"
$wsA.Range("F3").Value = "CS304"

$wsA.Range("B5").Value = "CS304"
$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "CS309"
$wsA.Range("E5").Value = "CS309"
$wsA.Range("F5").Value = "Free"

$wsA.Range("B6").Value = "Free"
$wsA.Range("E6").Value = "CS461"

$wsA.Range("C7").Value = "Free"
$wsA.Range("D7").Value = "CS303"
$wsA.Range("E7").Value = "Free"
$wsA.Range("F7").Value = "CS303"

# --- Section_B updates ---
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "Free"
$wsB.Range("C2").Value = "Free"
$wsB.Range("D2").Value = "CS461"
$wsB.Range("E2").Value = "CS461"
$wsB.Range("F2").Value = "CS303"

$wsB.Range("B3").Value = "Free"
$wsB.Range("C3").Value = "CS461"
$wsB.Range("D3").Value = "CS303"

$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "CS303"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "CS304"

$wsB.Range("B6").Value = "CS304"
$wsB.Range("D6").Value = "CS309"
$wsB.Range("F6").Value = "CS304"

$wsB.Range("B7").Value = "Free"
$wsB.Range("C7").Value = "CS309"
$wsB.Range("E7").Value = "Free"
